# Update database and change read_price algorithm:
# all previously-read numeric figures in the income-statement data block
# (rows 11-27, columns D:H) are reset - rows that used to hold a real
# numeric amount now read 0, while the rows/cells that model a "no value"
# placeholder (already shown as "-") now read "-" consistently across the
# whole row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "E", "F", "G", "H")

# Rows whose numeric figures are all zeroed out.
$zeroRows = @(11, 12, 13, 14, 16, 17, 19, 20, 21, 22, 24, 25, 26, 27)
foreach ($r in $zeroRows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = 0
    }
}

# Row 15: D15 was already "-"; E15:H15 become "-" too.
foreach ($c in @("E", "F", "G", "H")) {
    $ws.Range("${c}15").Value = "-"
}

# Row 18: D18 switches from 0 to "-"; E18:H18 are zeroed like the other rows.
$ws.Range("D18").Value = "-"
foreach ($c in @("E", "F", "G", "H")) {
    $ws.Range("${c}18").Value = 0
}

# Row 23: all of D23:H23 become "-".
foreach ($c in $cols) {
    $ws.Range("${c}23").Value = "-"
}
